# Update "want to go" counts (column F) on the "展览" (Exhibition) sheet
# and the "全部类型" (All types) sheet. Both sheets list the same set of
# exhibition events (全部类型 additionally contains a performance row),
# so the same F-column updates apply to each, just on different rows.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 27
    3  = 302
    5  = 2578
    6  = 1855
    8  = 109
    9  = 906
    10 = 179
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

# On "全部类型" the same events appear but shifted down by one row
# starting at row 8 (an extra performance entry occupies row 8 there).
$updates4 = @{
    2  = 27
    3  = 302
    5  = 2578
    6  = 1855
    9  = 109
    10 = 906
    11 = 179
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
